$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 319). Bump that serial value by one day (45188 -> 45189)
# for all of them, leaving everything else untouched.
$ws.Range("C2:C319").Value2 = 45189
